# Update "想去人数" (F column) values across the three affected worksheets.
# Sheet "本地生活" (local life) has no changes in this revision.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 636
$ws1.Range("F3").Value = 698
$ws1.Range("F4").Value = 943
$ws1.Range("F5").Value = 715
$ws1.Range("F7").Value = 400
$ws1.Range("F8").Value = 597
$ws1.Range("F9").Value = 129
$ws1.Range("F10").Value = 1204
$ws1.Range("F11").Value = 631
$ws1.Range("F12").Value = 382
$ws1.Range("F14").Value = 165
$ws1.Range("F15").Value = 10
$ws1.Range("F16").Value = 461
$ws1.Range("F17").Value = 352
$ws1.Range("F22").Value = 573
$ws1.Range("F23").Value = 27
$ws1.Range("F24").Value = 750

# --- Sheet: 演出 (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 85
$ws2.Range("F4").Value = 316
$ws2.Range("F5").Value = 101
$ws2.Range("F8").Value = 181
$ws2.Range("F9").Value = 221
$ws2.Range("F12").Value = 23
$ws2.Range("F13").Value = 93

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 85
$ws4.Range("F4").Value = 636
$ws4.Range("F6").Value = 316
$ws4.Range("F7").Value = 698
$ws4.Range("F8").Value = 943
$ws4.Range("F9").Value = 715
$ws4.Range("F11").Value = 400
$ws4.Range("F12").Value = 597
$ws4.Range("F13").Value = 129
$ws4.Range("F14").Value = 1204
$ws4.Range("F15").Value = 631
$ws4.Range("F16").Value = 101
$ws4.Range("F18").Value = 382
$ws4.Range("F21").Value = 165
$ws4.Range("F22").Value = 10
$ws4.Range("F23").Value = 461
$ws4.Range("F24").Value = 181
$ws4.Range("F25").Value = 352
$ws4.Range("F28").Value = 221
$ws4.Range("F32").Value = 23
$ws4.Range("F33").Value = 93
$ws4.Range("F35").Value = 573
$ws4.Range("F36").Value = 27
$ws4.Range("F37").Value = 750
